$d = $word.ActiveDocument

# Remove the leading empty paragraph.
$d.Paragraphs(1).Range.Delete()

# Remove the "Hello World!" Heading1 paragraph.
$d.Paragraphs(1).Range.Delete()

# The remaining paragraph ("New Paragraph " + several "This is so AwesomeN!"
# runs) becomes the Google Assistant response text. Replace its text
# (excluding the trailing paragraph mark) with the new sentence so the
# paragraph keeps a single run.
$lastPara = $d.Paragraphs(1)
$r = $lastPara.Range
$r.End = $r.End - 1
$r.Text = "megha is a bad girl"
